# Apply "Add data for 2021-11-12" update to the carjacking-by-neighborhood-by-month workbook.
# This advances the report's cutoff date from November 03 to November 04, 2021, and
# records newly-added carjacking incidents (current November 2021 plus matching-date
# incidents from prior Novembers) for several neighborhoods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the report title text (column B header / shared string).
$ws.Name = "Through 2021-11-04"
$ws.Range("B1").Value = "November 2021 (through November 04)"

# North Lawndale (row 2)
$ws.Range("B2").Value = 2
$ws.Range("AT2").Value = 1

# Austin (row 4)
$ws.Range("B4").Value = 1

# Englewood (row 7)
$ws.Range("X7").Value = 2

# Auburn Gresham (row 9)
$ws.Range("M9").Value = 1
$ws.Range("AT9").Value = 4

# Grand Boulevard (row 10)
$ws.Range("BE10").Value = 1

# New City (row 11)
$ws.Range("M11").Value = 1
$ws.Range("X11").Value = 1
$ws.Range("AT11").Value = 2

# Lower West Side (row 12)
$ws.Range("BE12").Value = 1

# Grand Crossing (row 15)
$ws.Range("B15").Value = 1

# Lake View (row 20)
$ws.Range("B20").Value = 1

# West Pullman (row 21)
$ws.Range("B21").Value = 1

# Albany Park (row 31)
$ws.Range("M31").Value = 1

# Lincoln Park (row 33)
$ws.Range("M33").Value = 3
$ws.Range("BP33").Value = 1

# Irving Park (row 34)
$ws.Range("AT34").Value = 1

# West Elsdon (row 39)
$ws.Range("BE39").Value = 1

# Fuller Park (row 58)
$ws.Range("BE58").Value = 1

# Chicago Lawn (row 66)
$ws.Range("M66").Value = 2

# Douglas (row 68)
$ws.Range("M68").Value = 2

# Jefferson Park (row 79)
$ws.Range("AI79").Value = 1

# Kenwood (row 80)
$ws.Range("B80").Value = 2

# Printers Row (row 90)
$ws.Range("AT90").Value = 1

# Rogers Park (row 91)
$ws.Range("AI91").Value = 1
